# "Enhance batch add certificate feature"
# - Rename/reorder columns: EXPO, JUDUL KEGIATAN, DURASI KEGIATAN, NAMA, ROLE, TANGGAL SERTIFIKAT
# - Add ROLE column (Peserta) and TANGGAL SERTIFIKAT column (2023-07-31)
# - Move the "(Baris ini hanya contoh, mohon dihapus)" note to column G
# - Re-style header row (center+middle) and the new date column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the existing example-row values before they get overwritten ---
$expo   = $ws.Range("A2").Value()
$judul  = $ws.Range("C2").Value()
$durasi = $ws.Range("D2").Value()
$nama   = $ws.Range("E2").Value()
$contoh = $ws.Range("F2").Value()

# --- header row (row 1) ---
$ws.Range("A1").Value = "EXPO"
$ws.Range("B1").Value = "JUDUL KEGIATAN"
$ws.Range("C1").Value = "DURASI KEGIATAN"
$ws.Range("D1").Value = "NAMA"
$ws.Range("E1").Value = "ROLE"
$ws.Range("F1").Value = "TANGGAL SERTIFIKAT"

# --- data row (row 2), re-ordered + two new columns ---
$ws.Range("A2").Value = $expo
$ws.Range("B2").Value = $judul
$ws.Range("C2").Value = $durasi
$ws.Range("D2").Value = $nama
$ws.Range("E2").Value = "Peserta"
$ws.Range("F2").Value = "2023-07-31"
$ws.Range("G2").Value = $contoh

# --- header styling: bold/fill already applied via s=2, add centered horizontal alignment ---
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").VerticalAlignment = -4108

# --- TANGGAL SERTIFIKAT header (F1): same look as other headers but text-format + wrap ---
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4108
$ws.Range("F1").WrapText = $true
$ws.Range("F1").NumberFormat = "@"

# --- TANGGAL SERTIFIKAT data cell (F2): text-format, centered ---
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").HorizontalAlignment = -4108

# --- note cell (G2) keeps the original italic style used at the old F2 ---
$ws.Range("G2").Font.Italic = $true

# --- column widths ---
$ws.Range("A1").ColumnWidth = 27.44140625
$ws.Range("B1").ColumnWidth = 36
$ws.Range("C1").ColumnWidth = 17.5546875
$ws.Range("D1").ColumnWidth = 24.88671875
$ws.Range("F1").ColumnWidth = 22.21875
$ws.Range("F1").EntireColumn.NumberFormat = "@"

# --- row 1 height ---
$ws.Range("A1:G1").RowHeight = 28.2

# --- selection / view niceties ---
$ws.Range("G6").Select()
